# Fruta / hortaliza, semanal
# Insert a new weekly record at row 319 (Vega Central Mapocho de Santiago - Mango),
# shifting the existing rows 319:337 down to 320:338.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("319:319").Insert()

$ws.Range("A319").Value = 9
$ws.Range("B319").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C319").Value = "Metropolitana"
$ws.Range("D319").Value = 44585
$ws.Range("E319").Value = 13
$ws.Range("F319").Value = "Fruta"
$ws.Range("G319").Value = 100108
$ws.Range("H319").Value = "Tropicales y subtropicales"
$ws.Range("I319").Value = 100108002
$ws.Range("J319").Value = "Mango"
$ws.Range("K319").Value = "Sin especificar"
$ws.Range("L319").Value = "Primera"
$ws.Range("M319").Value = 320
$ws.Range("N319").Value = 5500
$ws.Range("O319").Value = 6000
$ws.Range("P319").Value = 5781
$ws.Range("Q319").Value = '$/bandeja 4 kilos'
$ws.Range("R319").Value = "Perú"
$ws.Range("S319").Value = 1445
$ws.Range("T319").Value = 4
